# Weekly data refresh: insert the newest week's "Cebolla" price rows
# (Terminal La Palmera de La Serena) at the top of the historical block,
# pushing all existing rows for that block down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 1036-1037 (existing row 1036 and everything
# below shifts down to 1038 onward).
$ws.Range("A1036:A1037").EntireRow.Insert()

# New row 1036: Cebolla, Morada(o), 1a (cosecha) - $/malla 16 kilos
$ws.Range("A1036").Value = 8
$ws.Range("B1036").Value = "Terminal La Palmera de La Serena"
$ws.Range("C1036").Value = "Coquimbo"
$ws.Range("D1036").Value = 45041
$ws.Range("E1036").Value = 4
$ws.Range("F1036").Value = 100112004
$ws.Range("G1036").Value = "Cebolla"
$ws.Range("H1036").Value = "Morada(o)"
$ws.Range("I1036").Value = "1a (cosecha)"
$ws.Range("J1036").Value = 1600
$ws.Range("K1036").Value = 10000
$ws.Range("L1036").Value = 11000
$ws.Range("M1036").Value = 10500
$ws.Range("N1036").Value = "`$/malla 16 kilos"
$ws.Range("O1036").Value = "Región de O'Higgins"
$ws.Range("P1036").Value = 656
$ws.Range("Q1036").Value = 16
$ws.Range("R1036").Value = "Hortaliza"

# New row 1037: Cebolla, Sin especificar, 1a (cosecha) - $/malla 16 kilos
$ws.Range("A1037").Value = 8
$ws.Range("B1037").Value = "Terminal La Palmera de La Serena"
$ws.Range("C1037").Value = "Coquimbo"
$ws.Range("D1037").Value = 45041
$ws.Range("E1037").Value = 4
$ws.Range("F1037").Value = 100112004
$ws.Range("G1037").Value = "Cebolla"
$ws.Range("H1037").Value = "Sin especificar"
$ws.Range("I1037").Value = "1a (cosecha)"
$ws.Range("J1037").Value = 2000
$ws.Range("K1037").Value = 7800
$ws.Range("L1037").Value = 8000
$ws.Range("M1037").Value = 7900
$ws.Range("N1037").Value = "`$/malla 16 kilos"
$ws.Range("O1037").Value = "Región de O'Higgins"
$ws.Range("P1037").Value = 494
$ws.Range("Q1037").Value = 16
$ws.Range("R1037").Value = "Hortaliza"
